$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.647.66'
$ws.Range("E2").Value = '  -0.39%  '
$ws.Range("D3").Value = '3.446.10'
$ws.Range("E3").Value = '  -0.98%  '
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '1.00'
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = '  -0.02%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '579.64'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  -0.52%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '148.31'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +0.50%  '
$ws.Range("E7").Value = '  +0.03%  '
$ws.Range("E8").Value = '  +0.24%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '8.00'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +4.40%  '
$ws.Range("E10").Value = '  -1.56%  '
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("E13").Value = '  +1.95%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '28.32'
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = '  -5.45%  '
$ws.Range("E15").Value = '  -0.31%  '
$ws.Range("E16").Value = '  -0.16%  '
$ws.Range("D17").Value = '62.692.75'
$ws.Range("E17").Value = '  -0.36%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '6.33'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  -0.37%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '14.66'
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = '  +1.82%  '
$ws.Range("E20").Value = '  -2.95%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '386.83'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.86%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.564'
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = '  -0.66%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '75.38'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +0.23%  '
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").Value = '3.581.63'
$ws.Range("E25").Value = '  -1.20%  '
$ws.Range("E26").Value = '  -1.39%  '
$ws.Range("E27").Value = '  +1.22%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.66'
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = '  -0.33%  '
$ws.Range("E29").Value = '  +0.16%  '
$ws.Range("E30").Value = '  -2.91%  '
$ws.Range("E31").Value = '  -1.36%  '
$ws.Range("E33").Value = '  -4.62%  '
$ws.Range("E34").Value = '  -2.62%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.64'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +3.79%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '5.37'
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = '  +1.76%  '
$ws.Range("E37").Value = '  -0.06%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '6.95'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  -2.66%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '169.52'
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = '  -1.18%  '
$ws.Range("E40").Value = '  -1.01%  '
$ws.Range("E41").Value = '  +0.64%  '
$ws.Range("E42").Value = '  -2.38%  '
$ws.Range("E43").Value = '  +0.83%  '
$ws.Range("E44").Value = '  -1.51%  '
$ws.Range("E45").Value = '  -2.53%  '
$ws.Range("E46").Value = '  -1.70%  '
$ws.Range("D47").Value = '2.572.36'
$ws.Range("E47").Value = '  -1.47%  '
$ws.Range("E48").Value = '  +1.38%  '
$ws.Range("E49").Value = '  -1.08%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '22.70'
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = '  -3.78%  '
$ws.Range("E51").Value = '  +0.03%  '
